$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split "...sons pela suas respectivas " into three runs so the
# text reads "...sons pelas suas respectivas " (adds an "s" on its own run).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("pela suas respectivas frequências.") | Out-Null
$target1 = $d.Range($rng.Start, $rng.End)

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>pela</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> suas respectivas </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>fr</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>equências.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: rewrite the "Sinais Analógicos" section.
#   - merge the heading's two runs ("Sinais Analógicos" + ":") into one run
#   - replace the single explanatory paragraph + the three "%..." bullet
#     paragraphs + the trailing empty (bookmark) paragraph with five new
#     paragraphs describing analog signals.
# ---------------------------------------------------------------------------

# Heading paragraph: "Sinais Analógicos" + ":" -> "Sinais Analógicos:"
$headingPara = $d.Paragraphs(8)
$headingRange = $headingPara.Range
$headingTextRange = $d.Range($headingRange.Start, $headingRange.End - 1)

$xmlHeading = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:sz w:val="34"/><w:szCs w:val="34"/></w:rPr><w:t>Sinais Analógicos:</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$headingTextRange.InsertXML($xmlHeading)

# Body paragraphs 9-13 -> five new paragraphs.
$bodyStart = $d.Paragraphs(9).Range.Start
$bodyEnd = $d.Paragraphs(13).Range.End
$bodyRange = $d.Range($bodyStart, $bodyEnd)

$xmlBody = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Um sinal analógico </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">é uma onda variável </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">que representa uma quantidade variando em função do tempo, estes sinais normalmente </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">são utilizados em contexto elétrico, no entanto, também podem estar em um contexto mecânico, pneumático hidráulico e em muitos outros pois qualquer informação pode ser convertida em um sinal analógico. </w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Sinais analógicos podem ser usados para medir mudanças em fenômenos físicos como o som, a luz, a temperatura, a posição ou a pressão através de um transdutor de sinal, este têm basicamente a funcionalidade de converter energia de uma forma para outra. </w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A principal vantagem em se utilizar um sinal analógico é a boa definição deste sinal</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> pois ele possui uma quantidade infinita de resoluções, comparando com sinal digital percebemos que este possui uma maior densidade.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">A principal desvantagem em ser utilizar um sinal analógico é que este possui ruídos. </w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>%Base de todos os sistemas de telecomu</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>nicações;</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$bodyRange.InsertXML($xmlBody)

Write-Host "done"
